$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.633.09"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").Value = "3.309.05"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'572.00"
$ws.Range("E5").Value = "  -3.23%  "
$ws.Range("D6").Value = "'182.32"
$ws.Range("E6").Value = "  -5.39%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("D10").Value = "'6.62"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").Value = "'0.403"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("D12").Value = "3.888.22"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "'27.11"
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("D15").Value = "66.689.93"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "3.309.13"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "'13.73"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'432.75"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "'7.62"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Value = "'73.75"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'0.515"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").Value = "'9.06"
$ws.Range("E27").Value = "  -5.28%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "'22.81"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "'5.32"
$ws.Range("E31").Value = "  -4.52%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'6.77"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("D34").Value = "'1.22"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").Value = "'1.50"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").Value = "'160.32"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("D38").Value = "'27.30"
$ws.Range("D39").Value = "2.809.84"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "'0.790"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").Value = "'4.45"
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").Value = "'0.0675"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("D44").Value = "'40.13"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").Value = "'24.33"
$ws.Range("E45").Value = "  -3.59%  "
$ws.Range("D46").Value = "'2.34"
$ws.Range("E46").Value = "  -6.71%  "
$ws.Range("D47").Value = "'319.75"
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").Value = "'0.983"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").Value = "'6.17"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("E51").Value = "  -1.68%  "
